$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '41.521.69'
$ws.Range('E2').Value = '  +0.01%  '
$ws.Range('D3').Value = '2.465.71'
$ws.Range('E3').Value = '  -0.89%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.00'
$ws.Range('D4').ClearFormats()
$ws.Range('E4').Value = '  +0.02%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '314.43'
$ws.Range('D5').ClearFormats()
$ws.Range('E5').Value = '  +0.23%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '91.58'
$ws.Range('D6').ClearFormats()
$ws.Range('E6').Value = '  -3.63%  '
$ws.Range('E7').Value = '  -0.25%  '
$ws.Range('E8').Value = '  -0.06%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.510'
$ws.Range('D9').ClearFormats()
$ws.Range('E9').Value = '  +1.77%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '32.12'
$ws.Range('D10').ClearFormats()
$ws.Range('E10').Value = '  -4.66%  '
$ws.Range('E11').Value = '  +0.90%  '
$ws.Range('E12').Value = '  +0.55%  '
$ws.Range('D13').Value = '2.848.85'
$ws.Range('E13').Value = '  -0.76%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '6.82'
$ws.Range('D14').ClearFormats()
$ws.Range('E14').Value = '  -3.00%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '15.93'
$ws.Range('D15').ClearFormats()
$ws.Range('E15').Value = '  +3.01%  '
$ws.Range('D16').Value = '2.495.42'
$ws.Range('E16').Value = '  -0.24%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.769'
$ws.Range('D17').ClearFormats()
$ws.Range('E17').Value = '  -3.43%  '
$ws.Range('D18').Value = '41.497.93'
$ws.Range('E18').Value = '  +0.04%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '6.46'
$ws.Range('D19').ClearFormats()
$ws.Range('E19').Value = '  +1.70%  '
$ws.Range('D20').Value = '0.0₃0940'
$ws.Range('E20').Value = '  +1.58%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '71.05'
$ws.Range('D21').ClearFormats()
$ws.Range('E21').Value = '  +2.86%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '11.05'
$ws.Range('D22').ClearFormats()
$ws.Range('E22').Value = '  -2.53%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '235.28'
$ws.Range('D23').ClearFormats()
$ws.Range('E23').Value = '  -0.95%  '
$ws.Range('E24').Value = '  -1.95%  '
$ws.Range('E25').Value = '  -0.05%  '
$ws.Range('E26').Value = '  -0.92%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '24.55'
$ws.Range('D27').ClearFormats()
$ws.Range('E27').Value = '  +1.03%  '
$ws.Range('E28').Value = '  -0.41%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '9.64'
$ws.Range('D29').ClearFormats()
$ws.Range('E29').Value = '  -1.60%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '35.24'
$ws.Range('D30').ClearFormats()
$ws.Range('E30').Value = '  -3.98%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '155.55'
$ws.Range('D31').ClearFormats()
$ws.Range('E31').Value = '  +1.78%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '5.42'
$ws.Range('D32').ClearFormats()
$ws.Range('E32').Value = '  -1.90%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '2.56'
$ws.Range('D33').ClearFormats()
$ws.Range('E33').Value = '  -0.62%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.0757'
$ws.Range('D34').ClearFormats()
$ws.Range('E34').Value = '  -0.29%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '17.26'
$ws.Range('D35').ClearFormats()
$ws.Range('E35').Value = '  -5.33%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '2.86'
$ws.Range('D36').ClearFormats()
$ws.Range('E36').Value = '  -7.57%  '
$ws.Range('E37').Value = '  +1.09%  '
$ws.Range('E38').Value = '  -1.01%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '1.76'
$ws.Range('D39').ClearFormats()
$ws.Range('E39').Value = '  -6.60%  '
$ws.Range('E40').Value = '  -12.36%  '
$ws.Range('E41').Value = '  -4.42%  '
$ws.Range('E42').Value = '  -0.06%  '
$ws.Range('D43').Value = '1.938.10'
$ws.Range('E43').Value = '  -3.21%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.0282'
$ws.Range('D44').ClearFormats()
$ws.Range('E44').Value = '  -2.13%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '18.32'
$ws.Range('D45').ClearFormats()
$ws.Range('E46').Value = '  -4.30%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '9.01'
$ws.Range('D47').ClearFormats()
$ws.Range('E47').Value = '  +1.96%  '
$ws.Range('D48').Value = '2.708.52'
$ws.Range('E48').Value = '  -0.92%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '96.60'
$ws.Range('D49').ClearFormats()
$ws.Range('E49').Value = '  -0.90%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '66.85'
$ws.Range('D50').ClearFormats()
$ws.Range('E50').Value = '  -4.69%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '52.32'
$ws.Range('D51').ClearFormats()
$ws.Range('E51').Value = '  +2.24%  '
